# Update main GSC export data:
#  - drop the oldest date row (2025-10-20), which shifts every remaining
#    row's C (HTTPS URLs) value up by one date
#  - append 4 new trailing date rows (2026-01-16 .. 2026-01-19)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")

# Remove the first data row (row 2 = 2025-10-20); this shifts all rows
# below it up by one, so each surviving date keeps the next day's count.
$ws.Rows.Item(2).Delete()

# Figure out where the data now ends so we can append after it.
$used = $ws.UsedRange
$lastRow = $used.Rows.Count

$newDates = @("2026-01-16", "2026-01-17", "2026-01-18", "2026-01-19")
$newCounts = @(25, 25, 25, 25)

for ($i = 0; $i -lt $newDates.Length; $i++) {
    $r = $lastRow + 1 + $i

    $dateCell = $ws.Cells.Item($r, 1)
    # Prefix with an apostrophe so the date-looking string is kept as
    # plain text instead of being auto-parsed into a date serial value,
    # then strip the resulting formatting so the cell matches the plain
    # (unstyled) text cells used elsewhere in the column.
    $dateCell.Value = "'" + $newDates[$i]
    $dateCell.ClearFormats()

    $ws.Cells.Item($r, 2).Value = 0
    $ws.Cells.Item($r, 3).Value = $newCounts[$i]
}
